$d = $word.ActiveDocument

# Locate the paragraph containing "Creating a visual chart with numbers" in the
# "Predicting Fingers" problem-solving section, so the new bullet point can be
# inserted directly below it.
$anchor = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Creating a visual chart with numbers") {
        $anchor = $para
    }
}

if ($anchor -ne $null) {
    # Create a new paragraph right after the anchor; it inherits the anchor's
    # paragraph/character formatting (720-twip left indent, black/text1 color).
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($anchor.Index + 1)
    $newPara.Range.Text = [char]0x2022 + " Can multiply by either 10, 20, or even 50"
}
